$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 313, shifting existing rows 313:335 down to 314:336
$ws.Rows(313).Insert()

# Populate the newly inserted row 313 with the new weekly data point.
# Columns A,B,C,E,F,G,H,I,K,N,O,Q,R keep the same constant values as the
# rest of this Ciboulette / Femacal de La Calera dataset.
$ws.Range("A313").Value = 3
$ws.Range("B313").Value = "Femacal de La Calera"
$ws.Range("C313").Value = "Coquimbo"
$ws.Range("D313").Value = 44746
$ws.Range("E313").Value = 5
$ws.Range("F313").Value = 100112039
$ws.Range("G313").Value = "Ciboulette"
$ws.Range("H313").Value = "Sin especificar"
$ws.Range("I313").Value = "Primera"
$ws.Range("J313").Value = 135
$ws.Range("K313").Value = 1500
$ws.Range("L313").Value = 1800
$ws.Range("M313").Value = 1667
$ws.Range("N313").Value = "`$/docena de atados"
$ws.Range("O313").Value = "Provincia de Quillota"
$ws.Range("P313").Value = 556
$ws.Range("Q313").Value = 3
$ws.Range("R313").Value = "Hortaliza"
